$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains its text formatting so numeric-looking
# strings (e.g. "301.76", "1.002") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.894.20"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").Value = "1.875.02"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "301.76"
$ws.Range("E5").Value = "  -1.76%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "0.5323"
$ws.Range("E7").Value = "  +1.78%  "
$ws.Range("D8").Value = "0.3756"
$ws.Range("E8").Value = "  -1.23%  "
$ws.Range("D9").Value = "0.07166"
$ws.Range("E9").Value = "  -1.79%  "
$ws.Range("D10").Value = "21.60"
$ws.Range("E10").Value = "  +1.24%  "
$ws.Range("D11").Value = "0.8858"
$ws.Range("E11").Value = "  -2.23%  "
$ws.Range("D12").Value = "0.08135"
$ws.Range("E12").Value = "  -0.77%  "
$ws.Range("D13").Value = "1.857.89"
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("D14").Value = "93.15"
$ws.Range("E14").Value = "  -2.19%  "
$ws.Range("D15").Value = "5.276"
$ws.Range("E15").Value = "  -1.29%  "
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "14.75"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").Value = "0.000008546"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("D20").Value = "27.088.00"
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").Value = "4.974"
$ws.Range("E21").Value = "  -2.75%  "
$ws.Range("D22").Value = "10.70"
$ws.Range("E22").Value = "  -0.81%  "
$ws.Range("D23").Value = "6.402"
$ws.Range("E23").Value = "  -1.07%  "
$ws.Range("D24").Value = "147.49"
$ws.Range("E24").Value = "  -1.41%  "
$ws.Range("D25").Value = "2.276"
$ws.Range("E25").Value = "  -2.77%  "
$ws.Range("D26").Value = "1.743"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "18.02"
$ws.Range("E27").Value = "  -1.32%  "
$ws.Range("D28").Value = "114.64"
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("D29").Value = "4.742"
$ws.Range("E29").Value = "  -1.62%  "
$ws.Range("D30").Value = "4.596"
$ws.Range("E30").Value = "  -5.32%  "
$ws.Range("D31").Value = "0.09102"
$ws.Range("E31").Value = "  -1.48%  "
$ws.Range("D32").Value = "0.8010"
$ws.Range("E32").Value = "  +0.99%  "
$ws.Range("D33").Value = "0.04989"
$ws.Range("E33").Value = "  -1.21%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "2.995"
$ws.Range("E34").Value = "  +1.25%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "1.174"
$ws.Range("E35").Value = "  -3.90%  "
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").Value = "3.210"
$ws.Range("E36").Value = "  -5.24%  "
$ws.Range("B37").Value = "TheSandbox"
$ws.Range("C37").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D37").Value = "0.5851"
$ws.Range("E37").Value = "  +2.19%  "
$ws.Range("D38").Value = "2.592"
$ws.Range("E38").Value = "  -2.62%  "
$ws.Range("D39").Value = "0.01956"
$ws.Range("E39").Value = "  -1.94%  "
$ws.Range("D40").Value = "1.068"
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "6.599"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "8.928"
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("D43").Value = "116.51"
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("D44").Value = "0.5043"
$ws.Range("E44").Value = "  +3.01%  "
$ws.Range("D45").Value = "0.1496"
$ws.Range("E45").Value = "  -1.38%  "
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").Value = "9.965"
$ws.Range("E47").Value = "  -2.04%  "
$ws.Range("D48").Value = "1.610"
$ws.Range("D49").Value = "37.93"
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("D50").Value = "0.06030"
$ws.Range("E50").Value = "  +1.20%  "
$ws.Range("E51").Value = "  -2.32%  "
